$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.296.12"
$ws.Range("E2").Value = "  +3.54%  "

$ws.Range("D3").Value = "3.193.79"
$ws.Range("E3").Value = "  +5.10%  "

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.Value = "'0.999"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  -0.05%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.Value = "'205.65"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +1.81%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.Value = "'635.14"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +0.40%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.Value = "'0.999"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -0.01%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.Value = "'0.230"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +9.44%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.Value = "'0.584"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +5.30%  "

$ws.Range("D10").Value = "3.189.58"
$ws.Range("E10").Value = "  +5.13%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.Value = "'0.577"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +31.62%  "

$ws.Range("E12").Value = "  +3.04%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.Value = "'5.51"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +6.16%  "

$ws.Range("D14").Value = "3.778.38"
$ws.Range("E14").Value = "  +5.22%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.Value = "'0.0000226"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +15.42%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.Value = "'31.69"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +7.00%  "

$ws.Range("D17").Value = "79.102.68"
$ws.Range("E17").Value = "  +3.41%  "

$ws.Range("D18").Value = "3.185.30"
$ws.Range("E18").Value = "  +5.46%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.Value = "'14.50"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +7.39%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.Value = "'3.07"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  +32.48%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.Value = "'9.14"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +0.95%  "

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.Value = "'427.78"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +13.48%  "

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.Value = "'5.00"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +13.95%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.Value = "'6.87"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +6.11%  "

$ws.Range("E25").Value = "  +8.89%  "

$ws.Range("D26").Value = "3.358.47"
$ws.Range("E26").Value = "  +5.48%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.Value = "'11.19"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  +11.80%  "

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.Value = "'76.58"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  +3.69%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.Value = "'1.00"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +0.02%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.Value = "'0.0000117"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +2.18%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.Value = "'0.995"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  +7.65%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.Value = "'1.50"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  +5.12%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.Value = "'522.66"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("E35").Value = "  +1.58%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.Value = "'0.140"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  +23.92%  "

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.Value = "'22.90"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +9.49%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.Value = "'0.121"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +12.93%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("E40").Value = "  +4.18%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.Value = "'164.42"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  +0.88%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.Value = "'20.01"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -0.06%  "

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.Value = "'192.70"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +2.58%  "

$ws.Range("E44").Value = "  -0.20%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.Value = "'5.47"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +5.50%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.Value = "'0.812"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +11.44%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.Value = "'1.80"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  +6.89%  "

$ws.Range("E48").Value = "  +4.74%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.Value = "'42.73"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +1.21%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.Value = "'25.86"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +13.68%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.Value = "'2.52"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +1.44%  "
